$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H38").Value = 25117.428
$ws.Range("I38").Value = 274.33334
$ws.Range("J38").Value = 43749.75
$ws.Range("K38").Value = 823.0000200000001
$ws.Range("L38").Value = 131249.25
$ws.Range("M38").Value = -451.0000200000001
$ws.Range("N38").Value = -131993.25
$ws.Range("H98").Value = 2774.7837
$ws.Range("I98").Value = 2813.5588
$ws.Range("K98").Value = 2813.5588
$ws.Range("M98").Value = -1315.5588
$ws.Range("H113").Value = 2498.75
$ws.Range("I113").Value = 1999
$ws.Range("J113").Value = 3998
$ws.Range("K113").Value = 1999
$ws.Range("L113").Value = 3998
$ws.Range("M113").Value = 1255
$ws.Range("N113").Value = -10506
$ws.Range("H122").Value = 2774.7837
$ws.Range("I122").Value = 2813.5588
$ws.Range("K122").Value = 8440.6764
$ws.Range("M122").Value = -5990.6764
$ws.Range("H126").Value = 77737.5
$ws.Range("J126").Value = 77737.5
$ws.Range("L126").Value = 77737.5
$ws.Range("N126").Value = -87617.5
$ws.Range("H137").Value = 2250.3215
$ws.Range("I137").Value = 2050.138
$ws.Range("J137").Value = 2465.3333
$ws.Range("K137").Value = 6150.414
$ws.Range("L137").Value = 7395.999899999999
$ws.Range("M137").Value = -3600.414
$ws.Range("N137").Value = -12495.9999
$ws.Range("H138").Value = 2279914
$ws.Range("J138").Value = 3132632
$ws.Range("L138").Value = 9397896
$ws.Range("N138").Value = -9408176

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13361.6045
$ws.Range("I32").Value = 12039.743
$ws.Range("J32").Value = 26249.75
$ws.Range("K32").Value = 12039.743
$ws.Range("L32").Value = 26249.75
$ws.Range("M32").Value = -11752.743
$ws.Range("N32").Value = -26823.75
$ws.Range("H63").Value = 3096.9375
$ws.Range("I63").Value = 2242.4614
$ws.Range("J63").Value = 6799.6665
$ws.Range("K63").Value = 2242.4614
$ws.Range("L63").Value = 6799.6665
$ws.Range("M63").Value = -1556.4614
$ws.Range("N63").Value = -8171.6665
$ws.Range("H66").Value = 3096.9375
$ws.Range("I66").Value = 2242.4614
$ws.Range("J66").Value = 6799.6665
$ws.Range("K66").Value = 11212.307
$ws.Range("L66").Value = 33998.3325
$ws.Range("M66").Value = -7780.307000000001
$ws.Range("N66").Value = -40862.3325
$ws.Range("H74").Value = 220517.66
$ws.Range("I74").Value = 300339.94
$ws.Range("K74").Value = 300339.94
$ws.Range("M74").Value = -299465.94
$ws.Range("H77").Value = 220517.66
$ws.Range("I77").Value = 300339.94
$ws.Range("K77").Value = 1501699.7
$ws.Range("M77").Value = -1497331.7
$ws.Range("H97").Value = 2481.75
$ws.Range("I97").Value = 2481.75
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 2481.75
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1985.75
$ws.Range("N97").ClearContents()
$ws.Range("H132").Value = 3491.1538
$ws.Range("I132").Value = 3006.5881
$ws.Range("J132").Value = 4406.4443
$ws.Range("K132").Value = 9019.764299999999
$ws.Range("L132").Value = 13219.3329
$ws.Range("M132").Value = -6489.764299999999
$ws.Range("N132").Value = -18279.3329

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 50000
$ws.Range("J16").Value = 50000
$ws.Range("L16").Value = 50000
$ws.Range("N16").Value = -50340
$ws.Range("H20").Value = 37883670
$ws.Range("I20").Value = 49025336
$ws.Range("J20").Value = 2013.8
$ws.Range("K20").Value = 49025336
$ws.Range("L20").Value = 2013.8
$ws.Range("M20").Value = -49025089
$ws.Range("N20").Value = -2507.8
$ws.Range("H105").Value = 43336616
$ws.Range("I105").Value = 10000000
$ws.Range("K105").Value = 10000000
$ws.Range("M105").Value = -9998253

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 11632103
$ws.Range("I132").Value = 13516942
$ws.Range("K132").Value = 40550826
$ws.Range("M132").Value = -40548296
$ws.Range("H134").Value = 6186.477
$ws.Range("I134").Value = 6535.615
$ws.Range("J134").Value = 5682.1665
$ws.Range("K134").Value = 19606.845
$ws.Range("L134").Value = 17046.4995
$ws.Range("M134").Value = -17071.845
$ws.Range("N134").Value = -22116.4995
$ws.Range("H140").Value = 69927.14
$ws.Range("J140").Value = 69998.46000000001
$ws.Range("L140").Value = 69998.46000000001
$ws.Range("N140").Value = -80358.46000000001
$ws.Range("H141").Value = 337520.78
$ws.Range("J141").Value = 337520.78
$ws.Range("L141").Value = 337520.78
$ws.Range("N141").Value = -347880.78

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 356
$ws.Range("I50").Value = 490
$ws.Range("J50").Value = 88
$ws.Range("K50").Value = 1470
$ws.Range("L50").Value = 264
$ws.Range("M50").Value = -989
$ws.Range("N50").Value = -1226
$ws.Range("H53").Value = 356
$ws.Range("I53").Value = 490
$ws.Range("J53").Value = 88
$ws.Range("K53").Value = 1470
$ws.Range("L53").Value = 264
$ws.Range("M53").Value = -989
$ws.Range("N53").Value = -1226
$ws.Range("H129").Value = 3130.6667
$ws.Range("J129").Value = 2492.647
$ws.Range("L129").Value = 7477.941
$ws.Range("N129").Value = -17477.941

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 125525000
$ws.Range("I70").Value = 167333330
$ws.Range("J70").Value = 100000
$ws.Range("K70").Value = 167333330
$ws.Range("L70").Value = 100000
$ws.Range("M70").Value = -167333060
$ws.Range("N70").Value = -100540
$ws.Range("H73").Value = 125525000
$ws.Range("I73").Value = 167333330
$ws.Range("J73").Value = 100000
$ws.Range("K73").Value = 167333330
$ws.Range("L73").Value = 100000
$ws.Range("M73").Value = -167332394
$ws.Range("N73").Value = -101872
$ws.Range("H102").Value = 2967.8
$ws.Range("I102").Value = 2963
$ws.Range("K102").Value = 2963
$ws.Range("M102").Value = -1341
$ws.Range("H104").Value = 28670.5
$ws.Range("J104").Value = 28670.5
$ws.Range("L104").Value = 28670.5
$ws.Range("N104").Value = -35658.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3891.12
$ws.Range("I132").Value = 2786.6875
$ws.Range("J132").Value = 5854.5557
$ws.Range("K132").Value = 8360.0625
$ws.Range("L132").Value = 17563.6671
$ws.Range("M132").Value = -5830.0625
$ws.Range("N132").Value = -22623.6671
$ws.Range("H136").Value = 6935.1924
$ws.Range("I136").Value = 4839.7827
$ws.Range("J136").Value = 23000
$ws.Range("K136").Value = 14519.3481
$ws.Range("L136").Value = 69000
$ws.Range("M136").Value = -11969.3481
$ws.Range("N136").Value = -74100
$ws.Range("H139").Value = 78830.125
$ws.Range("I139").Value = 65323
$ws.Range("J139").Value = 83332.5
$ws.Range("K139").Value = 65323
$ws.Range("L139").Value = 83332.5
$ws.Range("M139").Value = -60183
$ws.Range("N139").Value = -93612.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 10101010
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 10101010
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 10101010
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -10101234
$ws.Range("H51").Value = 33333
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H52").Value = 24999.875
$ws.Range("J52").Value = 25000
$ws.Range("L52").Value = 25000
$ws.Range("N52").Value = -25452
$ws.Range("H81").Value = 5571.5
$ws.Range("I81").Value = 4424.5
$ws.Range("J81").Value = 6489.1
$ws.Range("K81").Value = 8849
$ws.Range("L81").Value = 12978.2
$ws.Range("M81").Value = -7788
$ws.Range("N81").Value = -15100.2
$ws.Range("H84").Value = 5571.5
$ws.Range("I84").Value = 4424.5
$ws.Range("J84").Value = 6489.1
$ws.Range("K84").Value = 44245
$ws.Range("L84").Value = 64891
$ws.Range("M84").Value = -38941
$ws.Range("N84").Value = -75499
$ws.Range("H126").Value = 2718.889
$ws.Range("I126").Value = 2542.3333
$ws.Range("K126").Value = 7626.999899999999
$ws.Range("M126").Value = -5156.999899999999
$ws.Range("H132").Value = 4569810
$ws.Range("I132").Value = 5750768.5
$ws.Range("J132").Value = 3438.6
$ws.Range("K132").Value = 17252305.5
$ws.Range("L132").Value = 10315.8
$ws.Range("M132").Value = -17249775.5
$ws.Range("N132").Value = -15375.8
$ws.Range("H135").Value = 57410.25
$ws.Range("J135").Value = 57410.25
$ws.Range("L135").Value = 57410.25
$ws.Range("N135").Value = -67550.25
